# Applies the "Initial demonstration of AST structure" edit:
#  1. Adds a trailing space run to the Parsing body paragraph, then inserts
#     a brand-new paragraph about the grammar (carrying the _GoBack
#     bookmark) right after it, ahead of the existing blank paragraph.
#  2. Removes the _GoBack bookmark from the Weeding body paragraph (it
#     moved to the new paragraph above).
#  3. Appends a blank paragraph plus a new "AST" Heading1 section (five
#     body paragraphs) at the end of the document.

$d = $word.ActiveDocument

# --- Change 1: "Parsing" section -----------------------------------------
# Find the paragraph that ends the Parsing section body text.
$parsingBody = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*generated by the scanner.*") {
        $parsingBody = $p
    }
}

# Append a trailing space as its own run.
$parsingBody.Range.InsertAfter(" ")

# Insert a brand-new (blank) paragraph ahead of the paragraph that
# currently follows the Parsing body (the pre-existing empty paragraph),
# so the empty paragraph is preserved further down, then fill the new
# paragraph with the grammar text + the relocated bookmark.
$nextPara = $parsingBody.Next()
$nextPara.Range.InsertParagraphBefore()
$grammarPara = $parsingBody.Next()

$grammarXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t>The grammar was developed with extensive reference to the online Java documentation as well as the specifi</w:t></w:r>
<w:r><w:t>c details of the Joos language.</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
[void]$grammarPara.Range.InsertXML($grammarXml)

# --- Change 2: "Weeding" section -----------------------------------------
# Rewrite the Weeding body paragraph without the bookmark (it now lives in
# the new paragraph above), keeping text/runs identical.
$weedingBody = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*ensure that the program behaves in a legal manner.*") {
        $weedingBody = $p
    }
}

$weedingXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="002827E6" w:rsidRDefault="00C5682C" w:rsidP="002827E6">
<w:r><w:t xml:space="preserve">The weeder is designed around a set of individual weeds that each need to be checked. Every parse tree generated by the compiler is checked </w:t></w:r>
<w:r w:rsidR="00E23CB7"><w:t>by the weeder for</w:t></w:r>
<w:r w:rsidR="00F33706"><w:t xml:space="preserve"> each of the weeds to ensure tha</w:t></w:r>
<w:r w:rsidR="00E23CB7"><w:t>t the program behaves in a legal manner.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
[void]$weedingBody.Range.InsertXML($weedingXml)

# --- Change 3: new "AST" section at the end of the document --------------
$endRng = $d.Range($d.Content.End, $d.Content.End)

$astXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p/>
<w:p>
<w:pPr><w:pStyle w:val="Heading1"/></w:pPr>
<w:r><w:t>AST</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>The compiler will make use of an AST to simplify the parse tree generated by the parser. The AST is still in development, but the current design involves creating a class for each nonterminal in the AST.</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>These classes will have a meaningful name that should improve code readability for the rest of the parser.</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>Each class will contain a parent pointer and several child pointers. The class will have a child pointer for every possible type of child that that class could have. If a child is missing from a particular object, the pointer for that child will be initialized to NULL.</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>There will also be an epsilon class to represent non-terminals that are reduced to null.</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>The AST classes will eventually be tagged with attributes which can be used by attribute grammars throughout the remainder of the compiler.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
[void]$endRng.InsertXML($astXml)

Write-Output "Done."
